$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.0787
$ws.Range("C8").Value = -10.95579999999999
$ws.Range("C10").Value = -12.5792
$ws.Range("C12").Value = -14.3779
$ws.Range("D13").Value = -7.978300000000002
$ws.Range("C18").Value = -14.3096
$ws.Range("E20").Value = 12.97929999999999
$ws.Range("C25").Value = -10.96049999999999
